$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 36984.78942497371
$ws.Range("C2").Value = 1863.713111023834
$ws.Range("D2").Value = 5864.194604207496
$ws.Range("E2").Value = 2291.779877781744
$ws.Range("F2").Value = 221.0521250706038
$ws.Range("G2").Value = 799.587147699983
$ws.Range("H2").Value = 40.3001736459458
$ws.Range("I2").Value = 3091.367025481727
$ws.Range("J2").Value = 261.3522987165496
$ws.Range("K2").Value = 393.013095416546
$ws.Range("L2").Value = 21.91741730592575
$ws.Range("M2").Value = 278.0250160352341
$ws.Range("N2").Value = 41.65753511656836
$ws.Range("O2").Value = 35318.30623147275
$ws.Range("P2").Value = 2055.524019744857
$ws.Range("Q2").Value = 7035.737028209332
$ws.Range("R2").Value = 2898.498343348916
$ws.Range("S2").Value = 262.4617143176719
$ws.Range("T2").Value = 757.7201179376283
$ws.Range("U2").Value = 39.44895443994156
$ws.Range("V2").Value = 3656.218461286545
$ws.Range("W2").Value = 301.9106687576135
$ws.Range("X2").Value = 416.4765382603131
$ws.Range("Y2").Value = 21.35569289610375
$ws.Range("Z2").Value = 294.2538785565282
$ws.Range("AA2").Value = 55.32037164341867

$ws.Range("B3").Value = 35139.82353109887
$ws.Range("C3").Value = 1690.440770990897
$ws.Range("D3").Value = 5511.485981514406
$ws.Range("E3").Value = 2238.681337601321
$ws.Range("F3").Value = 214.2489973357439
$ws.Range("G3").Value = 840.5826305662382
$ws.Range("H3").Value = 40.02188985175228
$ws.Range("I3").Value = 3079.263968167559
$ws.Range("J3").Value = 254.2708871874962
$ws.Range("K3").Value = 362.3264014831789
$ws.Range("L3").Value = 24.12555135755219
$ws.Range("M3").Value = 260.5781337760715
$ws.Range("N3").Value = 41.33524022861447
$ws.Range("O3").Value = 37278.73176481328
$ws.Range("P3").Value = 1741.475776568035
$ws.Range("Q3").Value = 6918.282226475702
$ws.Range("R3").Value = 2708.275818993478
$ws.Range("S3").Value = 253.4574588476533
$ws.Range("T3").Value = 697.6124669116859
$ws.Range("U3").Value = 37.08361998645271
$ws.Range("V3").Value = 3405.888285905164
$ws.Range("W3").Value = 290.541078834106
$ws.Range("X3").Value = 444.7572931794446
$ws.Range("Y3").Value = 23.50821093050978
$ws.Range("Z3").Value = 328.8671087604378
$ws.Range("AA3").Value = 56.52843416597797

$ws.Range("B4").Value = 35670.40120009604
$ws.Range("C4").Value = 1748.937103663767
$ws.Range("D4").Value = 5243.529566463166
$ws.Range("E4").Value = 2044.955379152391
$ws.Range("F4").Value = 196.9924026273008
$ws.Range("G4").Value = 772.8018779050992
$ws.Range("H4").Value = 38.19431627864717
$ws.Range("I4").Value = 2817.75725705749
$ws.Range("J4").Value = 235.186718905948
$ws.Range("K4").Value = 369.251856386387
$ws.Range("L4").Value = 20.56326026971866
$ws.Range("M4").Value = 260.2173766191106
$ws.Range("N4").Value = 40.28671915993478
$ws.Range("O4").Value = 34967.84365323477
$ws.Range("P4").Value = 1955.997191947915
$ws.Range("Q4").Value = 6735.950429127836
$ws.Range("R4").Value = 2512.264471963871
$ws.Range("S4").Value = 227.7493592698368
$ws.Range("T4").Value = 654.4873509197188
$ws.Range("U4").Value = 34.32881547899326
$ws.Range("V4").Value = 3166.75182288359
$ws.Range("W4").Value = 262.0781747488301
$ws.Range("X4").Value = 396.2543559713862
$ws.Range("Y4").Value = 19.01366400801037
$ws.Range("Z4").Value = 284.1242905108756
$ws.Range("AA4").Value = 47.46093844235369

$ws.Range("B5").Value = 32647.87026909309
$ws.Range("C5").Value = 1763.767837903084
$ws.Range("D5").Value = 5403.07065083388
$ws.Range("E5").Value = 2028.882743576388
$ws.Range("F5").Value = 188.9873199756406
$ws.Range("G5").Value = 713.0034047802379
$ws.Range("H5").Value = 36.12928511257346
$ws.Range("I5").Value = 2741.886148356625
$ws.Range("J5").Value = 225.1166050882141
$ws.Range("K5").Value = 325.3385808165229
$ws.Range("L5").Value = 18.66445279530906
$ws.Range("M5").Value = 233.8347715843163
$ws.Range("N5").Value = 37.69567533837247
$ws.Range("O5").Value = 36186.72569018874
$ws.Range("P5").Value = 1884.033887926297
$ws.Range("Q5").Value = 7116.883714345494
$ws.Range("R5").Value = 2726.646221664719
$ws.Range("S5").Value = 253.1635916720752
$ws.Range("T5").Value = 713.6647968605633
$ws.Range("U5").Value = 37.89107176172122
$ws.Range("V5").Value = 3440.311018525283
$ws.Range("W5").Value = 291.0546634337964
$ws.Range("X5").Value = 410.6698631995355
$ws.Range("Y5").Value = 20.4323344750366
$ws.Range("Z5").Value = 300.0259382531694
$ws.Range("AA5").Value = 50.74899268971261

$ws.Range("B6").Value = 36974.38135238089
$ws.Range("C6").Value = 1971.081921566278
$ws.Range("D6").Value = 6237.721800697452
$ws.Range("E6").Value = 2404.758738272372
$ws.Range("F6").Value = 225.6470418362181
$ws.Range("G6").Value = 783.1023279153544
$ws.Range("H6").Value = 42.10090083547875
$ws.Range("I6").Value = 3187.861066187726
$ws.Range("J6").Value = 267.7479426716968
$ws.Range("K6").Value = 387.5376182218013
$ws.Range("L6").Value = 22.63332159269888
$ws.Range("M6").Value = 272.0246582558757
$ws.Range("N6").Value = 40.65699418550835
$ws.Range("O6").Value = 34037.08639854103
$ws.Range("P6").Value = 1969.814407356367
$ws.Range("Q6").Value = 6461.814865793849
$ws.Range("R6").Value = 2682.693670934113
$ws.Range("S6").Value = 245.1037549219388
$ws.Range("T6").Value = 715.0841630708672
$ws.Range("U6").Value = 37.05953411533455
$ws.Range("V6").Value = 3397.777834004981
$ws.Range("W6").Value = 282.1632890372734
$ws.Range("X6").Value = 398.4998238482113
$ws.Range("Y6").Value = 20.63306160534848
$ws.Range("Z6").Value = 271.3609392051098
$ws.Range("AA6").Value = 48.55508678535421

$ws.Range("B7").Value = 34048.34512309541
$ws.Range("C7").Value = 1809.977375102943
$ws.Range("D7").Value = 5533.32353787244
$ws.Range("E7").Value = 2159.542732544159
$ws.Range("F7").Value = 205.3246991204325
$ws.Range("G7").Value = 724.6160456364601
$ws.Range("H7").Value = 38.30731764814891
$ws.Range("I7").Value = 2884.158778180619
$ws.Range("J7").Value = 243.6320167685814
$ws.Range("K7").Value = 353.7184562348427
$ws.Range("L7").Value = 17.57795175567649
$ws.Range("M7").Value = 244.2240809511906
$ws.Range("N7").Value = 35.89492816239859
$ws.Range("O7").Value = 39748.91055695822
$ws.Range("P7").Value = 1954.231772442519
$ws.Range("Q7").Value = 7594.832236369612
$ws.Range("R7").Value = 2812.131488379547
$ws.Range("S7").Value = 259.3236854680563
$ws.Range("T7").Value = 711.5460650218032
$ws.Range("U7").Value = 38.13584943155273
$ws.Range("V7").Value = 3523.67755340135
$ws.Range("W7").Value = 297.459534899609
$ws.Range("X7").Value = 443.638665347588
$ws.Range("Y7").Value = 20.92133440938671
$ws.Range("Z7").Value = 330.0160218141575
$ws.Range("AA7").Value = 53.17140072714147

$ws.Range("B8").Value = 38666.44617415352
$ws.Range("C8").Value = 1761.108252905832
$ws.Range("D8").Value = 5907.780197285741
$ws.Range("E8").Value = 2453.299623878735
$ws.Range("F8").Value = 231.9725186381782
$ws.Range("G8").Value = 839.1770725630356
$ws.Range("H8").Value = 44.43655945661627
$ws.Range("I8").Value = 3292.476696441771
$ws.Range("J8").Value = 276.4090780947945
$ws.Range("K8").Value = 403.1707345705044
$ws.Range("L8").Value = 23.22314552223716
$ws.Range("M8").Value = 287.229812868405
$ws.Range("N8").Value = 44.94408914964191
$ws.Range("O8").Value = 36027.20407695029
$ws.Range("P8").Value = 1834.103387848691
$ws.Range("Q8").Value = 6243.634696270775
$ws.Range("R8").Value = 2353.072603634937
$ws.Range("S8").Value = 231.5665083891014
$ws.Range("T8").Value = 667.878425248301
$ws.Range("U8").Value = 35.06300640891956
$ws.Range("V8").Value = 3020.951028883238
$ws.Range("W8").Value = 266.6295147980209
$ws.Range("X8").Value = 413.1113148006331
$ws.Range("Y8").Value = 19.79208500210665
$ws.Range("Z8").Value = 297.0488495061834
$ws.Range("AA8").Value = 50.37113890699165

$ws.Range("B9").Value = 36440.72858284784
$ws.Range("C9").Value = 1833.717587287561
$ws.Range("D9").Value = 5958.15467724795
$ws.Range("E9").Value = 2324.180100690997
$ws.Range("F9").Value = 213.7833744146907
$ws.Range("G9").Value = 759.8475574042693
$ws.Range("H9").Value = 39.15836666954452
$ws.Range("I9").Value = 3084.027658095266
$ws.Range("J9").Value = 252.9417410842352
$ws.Range("K9").Value = 372.9624990529351
$ws.Range("L9").Value = 22.9772389903973
$ws.Range("M9").Value = 258.6483012491335
$ws.Range("N9").Value = 36.90503339355683
$ws.Range("O9").Value = 32822.24989419879
$ws.Range("P9").Value = 1972.533069187244
$ws.Range("Q9").Value = 6602.999732912004
$ws.Range("R9").Value = 2669.182464448235
$ws.Range("S9").Value = 235.132952968134
$ws.Range("T9").Value = 710.5611556902799
$ws.Range("U9").Value = 36.91028573641209
$ws.Range("V9").Value = 3379.743620138515
$ws.Range("W9").Value = 272.043238704546
$ws.Range("X9").Value = 360.4221816733501
$ws.Range("Y9").Value = 19.16336104346119
$ws.Range("Z9").Value = 250.370317805314
$ws.Range("AA9").Value = 44.42988666735041

$ws.Range("B10").Value = 38255.64517666001
$ws.Range("C10").Value = 1766.883441174571
$ws.Range("D10").Value = 6193.01590876589
$ws.Range("E10").Value = 2372.964735573008
$ws.Range("F10").Value = 224.3425409093067
$ws.Range("G10").Value = 827.9870874144821
$ws.Range("H10").Value = 41.58784447918433
$ws.Range("I10").Value = 3200.95182298749
$ws.Range("J10").Value = 265.930385388491
$ws.Range("K10").Value = 400.3020953679732
$ws.Range("L10").Value = 26.38165588884013
$ws.Range("M10").Value = 289.7676027772649
$ws.Range("N10").Value = 44.89919475651774
$ws.Range("O10").Value = 31960.22600064877
$ws.Range("P10").Value = 1812.256080771003
$ws.Range("Q10").Value = 5857.257841365194
$ws.Range("R10").Value = 2240.149115371444
$ws.Range("S10").Value = 206.092053379901
$ws.Range("T10").Value = 585.4340133463652
$ws.Range("U10").Value = 29.93613709859859
$ws.Range("V10").Value = 2825.583128717809
$ws.Range("W10").Value = 236.0281904784996
$ws.Range("X10").Value = 365.933849250119
$ws.Range("Y10").Value = 19.43182224566973
$ws.Range("Z10").Value = 265.413380104986
$ws.Range("AA10").Value = 44.40391729390301

$ws.Range("B11").Value = 34205.07958177093
$ws.Range("C11").Value = 1752.370680344522
$ws.Range("D11").Value = 5453.99385887174
$ws.Range("E11").Value = 2133.668123212671
$ws.Range("F11").Value = 202.9644532584234
$ws.Range("G11").Value = 744.842622276547
$ws.Range("H11").Value = 37.46209824611816
$ws.Range("I11").Value = 2878.510745489218
$ws.Range("J11").Value = 240.4265515045415
$ws.Range("K11").Value = 349.6933136790033
$ws.Range("L11").Value = 20.35895531075379
$ws.Range("M11").Value = 252.169763229969
$ws.Range("N11").Value = 36.17941948018166
$ws.Range("O11").Value = 38494.82309156028
$ws.Range("P11").Value = 1917.931615868652
$ws.Range("Q11").Value = 7360.582111863511
$ws.Range("R11").Value = 2925.93209082244
$ws.Range("S11").Value = 274.6926342690621
$ws.Range("T11").Value = 792.0766859198609
$ws.Range("U11").Value = 41.91618988838238
$ws.Range("V11").Value = 3718.008776742302
$ws.Range("W11").Value = 316.6088241574445
$ws.Range("X11").Value = 446.6336069517446
$ws.Range("Y11").Value = 22.60223208019178
$ws.Range("Z11").Value = 329.091952054256
$ws.Range("AA11").Value = 57.49436285072487

$ws.Range("B12").Value = 32758.08195563007
$ws.Range("C12").Value = 1793.357040425116
$ws.Range("D12").Value = 5238.694765575687
$ws.Range("E12").Value = 2086.253184673554
$ws.Range("F12").Value = 199.2119480330693
$ws.Range("G12").Value = 791.565795236153
$ws.Range("H12").Value = 38.89003550448751
$ws.Range("I12").Value = 2877.818979909707
$ws.Range("J12").Value = 238.1019835375568
$ws.Range("K12").Value = 345.8938173378686
$ws.Range("L12").Value = 20.40748480673266
$ws.Range("M12").Value = 237.5930698026506
$ws.Range("N12").Value = 37.69317383878818
$ws.Range("O12").Value = 38599.67811345107
$ws.Range("P12").Value = 1932.322711398921
$ws.Range("Q12").Value = 6718.351363627067
$ws.Range("R12").Value = 2380.928069859604
$ws.Range("S12").Value = 233.8774139636755
$ws.Range("T12").Value = 647.6339143920569
$ws.Range("U12").Value = 34.44766538418413
$ws.Range("V12").Value = 3028.561984251661
$ws.Range("W12").Value = 268.3250793478596
$ws.Range("X12").Value = 441.1664114612844
$ws.Range("Y12").Value = 20.96319904245146
$ws.Range("Z12").Value = 328.5689656180538
$ws.Range("AA12").Value = 53.81616349185092

$ws.Range("B13").Value = 34953.74543420662
$ws.Range("C13").Value = 1866.032217015614
$ws.Range("D13").Value = 5656.712370198617
$ws.Range("E13").Value = 2394.085879409804
$ws.Range("F13").Value = 222.110488791928
$ws.Range("G13").Value = 856.2483521944293
$ws.Range("H13").Value = 43.45043130386292
$ws.Range("I13").Value = 3250.334231604233
$ws.Range("J13").Value = 265.5609200957909
$ws.Range("K13").Value = 367.3673110889612
$ws.Range("L13").Value = 22.13109474157912
$ws.Range("M13").Value = 258.7621329727833
$ws.Range("N13").Value = 40.62741098809929
$ws.Range("O13").Value = 37707.68500705942
$ws.Range("P13").Value = 1927.470901623121
$ws.Range("Q13").Value = 7011.104832049384
$ws.Range("R13").Value = 2696.729128388443
$ws.Range("S13").Value = 251.0553186819429
$ws.Range("T13").Value = 725.6114142936802
$ws.Range("U13").Value = 39.06582525829596
$ws.Range("V13").Value = 3422.340542682123
$ws.Range("W13").Value = 290.1211439402388
$ws.Range("X13").Value = 431.0926321078739
$ws.Range("Y13").Value = 21.21677259777582
$ws.Range("Z13").Value = 307.5498847283098
$ws.Range("AA13").Value = 53.45620293251407

$ws.Range("B14").Value = 32963.50632179116
$ws.Range("C14").Value = 1767.23706594256
$ws.Range("D14").Value = 4769.842932808218
$ws.Range("E14").Value = 1895.187098614456
$ws.Range("F14").Value = 184.84857931165
$ws.Range("G14").Value = 679.4621046825976
$ws.Range("H14").Value = 32.77406582441426
$ws.Range("I14").Value = 2574.649203297054
$ws.Range("J14").Value = 217.6226451360643
$ws.Range("K14").Value = 333.4816257848676
$ws.Range("L14").Value = 19.66103671024132
$ws.Range("M14").Value = 237.581875238145
$ws.Range("N14").Value = 35.65024436258229
$ws.Range("O14").Value = 37916.45322392388
$ws.Range("P14").Value = 1904.955064933057
$ws.Range("Q14").Value = 6502.18421103032
$ws.Range("R14").Value = 2245.90156479402
$ws.Range("S14").Value = 218.1427384878003
$ws.Range("T14").Value = 557.7131031458366
$ws.Range("U14").Value = 29.90705738630099
$ws.Range("V14").Value = 2803.614667939856
$ws.Range("W14").Value = 248.0497958741013
$ws.Range("X14").Value = 409.2760007511265
$ws.Range("Y14").Value = 19.85950173582542
$ws.Range("Z14").Value = 304.509415916294
$ws.Range("AA14").Value = 50.69267749269447

$ws.Range("B15").Value = 32220.02150676567
$ws.Range("C15").Value = 1830.580065384755
$ws.Range("D15").Value = 5301.239069384605
$ws.Range("E15").Value = 2132.011747436541
$ws.Range("F15").Value = 198.2071088969121
$ws.Range("G15").Value = 695.8568040917461
$ws.Range("H15").Value = 35.4765987436576
$ws.Range("I15").Value = 2827.868551528287
$ws.Range("J15").Value = 233.6837076405697
$ws.Range("K15").Value = 336.9844681024867
$ws.Range("L15").Value = 20.6634454395392
$ws.Range("M15").Value = 241.4789991936919
$ws.Range("N15").Value = 36.72691545228616
$ws.Range("O15").Value = 38633.96516341814
$ws.Range("P15").Value = 1969.789367164996
$ws.Range("Q15").Value = 6997.280355470257
$ws.Range("R15").Value = 2697.922713108198
$ws.Range("S15").Value = 250.4019201540983
$ws.Range("T15").Value = 733.9124385027154
$ws.Range("U15").Value = 37.36432920592028
$ws.Range("V15").Value = 3431.835151610913
$ws.Range("W15").Value = 287.7662493600186
$ws.Range("X15").Value = 425.14171812629
$ws.Range("Y15").Value = 20.73584994464252
$ws.Range("Z15").Value = 313.5033248193728
$ws.Range("AA15").Value = 49.95718077567489

$ws.Range("B16").Value = 30965.54675942479
$ws.Range("C16").Value = 1800.850341017
$ws.Range("D16").Value = 4916.24135140615
$ws.Range("E16").Value = 1966.027655306044
$ws.Range("F16").Value = 183.1214673779982
$ws.Range("G16").Value = 677.8628917114975
$ws.Range("H16").Value = 33.30107361485849
$ws.Range("I16").Value = 2643.890547017541
$ws.Range("J16").Value = 216.4225409928567
$ws.Range("K16").Value = 323.6035574796436
$ws.Range("L16").Value = 17.61052936080061
$ws.Range("M16").Value = 223.6140481905484
$ws.Range("N16").Value = 33.918905803945
$ws.Range("O16").Value = 38398.02693295938
$ws.Range("P16").Value = 1890.408695170702
$ws.Range("Q16").Value = 7041.233320611997
$ws.Range("R16").Value = 2613.129243165958
$ws.Range("S16").Value = 240.5277961084957
$ws.Range("T16").Value = 698.5232165099505
$ws.Range("U16").Value = 34.20702297399577
$ws.Range("V16").Value = 3311.652459675908
$ws.Range("W16").Value = 274.7348190824914
$ws.Range("X16").Value = 432.9397177833973
$ws.Range("Y16").Value = 19.23656628065979
$ws.Range("Z16").Value = 318.9259233922268
$ws.Range("AA16").Value = 48.19573460368308

$ws.Range("B17").Value = 32092.43569744579
$ws.Range("C17").Value = 1842.455613025959
$ws.Range("D17").Value = 5183.70059725008
$ws.Range("E17").Value = 2048.634412762716
$ws.Range("F17").Value = 193.8997873659992
$ws.Range("G17").Value = 705.9473503914388
$ws.Range("H17").Value = 36.12607865931245
$ws.Range("I17").Value = 2754.581763154155
$ws.Range("J17").Value = 230.0258660253117
$ws.Range("K17").Value = 334.2046841778367
$ws.Range("L17").Value = 17.21672867397698
$ws.Range("M17").Value = 228.8672944249722
$ws.Range("N17").Value = 34.66413181414439
$ws.Range("O17").Value = 40525.21044240975
$ws.Range("P17").Value = 1970.158275158401
$ws.Range("Q17").Value = 7900.865829780558
$ws.Range("R17").Value = 3001.926673098896
$ws.Range("S17").Value = 271.3109752725626
$ws.Range("T17").Value = 806.4426957447115
$ws.Range("U17").Value = 41.23337589593452
$ws.Range("V17").Value = 3808.369368843608
$ws.Range("W17").Value = 312.5443511684971
$ws.Range("X17").Value = 444.2224908853827
$ws.Range("Y17").Value = 19.77996690059289
$ws.Range("Z17").Value = 333.5524342211514
$ws.Range("AA17").Value = 58.48504078741644

$ws.Range("B18").Value = 32706.80719737224
$ws.Range("C18").Value = 1854.540548606246
$ws.Range("D18").Value = 5486.705109831577
$ws.Range("E18").Value = 2073.734208818265
$ws.Range("F18").Value = 203.242568382218
$ws.Range("G18").Value = 621.2318161859457
$ws.Range("H18").Value = 35.01734905756427
$ws.Range("I18").Value = 2694.96602500421
$ws.Range("J18").Value = 238.2599174397823
$ws.Range("K18").Value = 341.7040212931961
$ws.Range("L18").Value = 19.8393154251435
$ws.Range("M18").Value = 237.7130449073506
$ws.Range("N18").Value = 38.00438393296406
$ws.Range("O18").Value = 41842.81071672353
$ws.Range("P18").Value = 2147.247733113258
$ws.Range("Q18").Value = 9202.918048062398
$ws.Range("R18").Value = 3992.576644895034
$ws.Range("S18").Value = 352.9747693873791
$ws.Range("T18").Value = 1100.605823962163
$ws.Range("U18").Value = 60.41563664314325
$ws.Range("V18").Value = 5093.182468857197
$ws.Range("W18").Value = 413.3904060305224
$ws.Range("X18").Value = 467.8653562805059
$ws.Range("Y18").Value = 22.64911874439506
$ws.Range("Z18").Value = 325.7334128163661
$ws.Range("AA18").Value = 67.31286071329997

$ws.Range("B19").Value = 34994.11276809201
$ws.Range("C19").Value = 1886.173303597148
$ws.Range("D19").Value = 5831.629067891557
$ws.Range("E19").Value = 2289.377205619548
$ws.Range("F19").Value = 219.4289286512258
$ws.Range("G19").Value = 815.779541114235
$ws.Range("H19").Value = 42.00194898887105
$ws.Range("I19").Value = 3105.156746733783
$ws.Range("J19").Value = 261.4308776400969
$ws.Range("K19").Value = 364.9969406796904
$ws.Range("L19").Value = 22.33109916375548
$ws.Range("M19").Value = 251.4333672983258
$ws.Range("N19").Value = 38.02113080195703
$ws.Range("O19").Value = 35593.54238090029
$ws.Range("P19").Value = 1956.507447167411
$ws.Range("Q19").Value = 6581.101173174808
$ws.Range("R19").Value = 2583.964680104031
$ws.Range("S19").Value = 247.1156971093195
$ws.Range("T19").Value = 746.6306666501458
$ws.Range("U19").Value = 39.28784236278945
$ws.Range("V19").Value = 3330.595346754177
$ws.Range("W19").Value = 286.4035394721089
$ws.Range("X19").Value = 411.6386442342052
$ws.Range("Y19").Value = 20.93217797149137
$ws.Range("Z19").Value = 290.2335250147322
$ws.Range("AA19").Value = 49.70419316230832

$ws.Range("B20").Value = 33635.70042211857
$ws.Range("C20").Value = 1841.903180168291
$ws.Range("D20").Value = 5305.916087443142
$ws.Range("E20").Value = 2031.035600461399
$ws.Range("F20").Value = 191.4129089669059
$ws.Range("G20").Value = 655.4278370891473
$ws.Range("H20").Value = 33.78516259738231
$ws.Range("I20").Value = 2686.463437550546
$ws.Range("J20").Value = 225.1980715642882
$ws.Range("K20").Value = 354.8949109155481
$ws.Range("L20").Value = 18.95430460545614
$ws.Range("M20").Value = 248.556750185756
$ws.Range("N20").Value = 36.54961972291465
$ws.Range("O20").Value = 38194.24051650666
$ws.Range("P20").Value = 1887.986747529322
$ws.Range("Q20").Value = 6792.473984308117
$ws.Range("R20").Value = 2444.004356279804
$ws.Range("S20").Value = 232.3109276618534
$ws.Range("T20").Value = 640.7241277567756
$ws.Range("U20").Value = 33.65720504742838
$ws.Range("V20").Value = 3084.728484036579
$ws.Range("W20").Value = 265.9681327092818
$ws.Range("X20").Value = 421.0167316122428
$ws.Range("Y20").Value = 20.40228236563006
$ws.Range("Z20").Value = 312.3413544941554
$ws.Range("AA20").Value = 48.51663485661462

$ws.Range("B21").Value = 34289.36642145396
$ws.Range("C21").Value = 1785.135697777305
$ws.Range("D21").Value = 5425.660954276913
$ws.Range("E21").Value = 2109.576997424297
$ws.Range("F21").Value = 197.5925635397815
$ws.Range("G21").Value = 737.2301563109829
$ws.Range("H21").Value = 37.39092385780267
$ws.Range("I21").Value = 2846.80715373528
$ws.Range("J21").Value = 234.9834873975841
$ws.Range("K21").Value = 353.6356650982439
$ws.Range("L21").Value = 19.43139192956996
$ws.Range("M21").Value = 250.3691403750342
$ws.Range("N21").Value = 36.53557994697282
$ws.Range("O21").Value = 41663.34866694458
$ws.Range("P21").Value = 1923.054261902613
$ws.Range("Q21").Value = 7949.114128038998
$ws.Range("R21").Value = 2965.687650254197
$ws.Range("S21").Value = 277.0982982233472
$ws.Range("T21").Value = 781.4176108836199
$ws.Range("U21").Value = 41.08631266665343
$ws.Range("V21").Value = 3747.105261137817
$ws.Range("W21").Value = 318.1846108900007
$ws.Range("X21").Value = 462.5051145262245
$ws.Range("Y21").Value = 20.58648107071341
$ws.Range("Z21").Value = 352.7150321250502
$ws.Range("AA21").Value = 54.79992878442881
